# Update "想去人数" (interest count) figures in column F on the "展览"
# and "全部类型" worksheets, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

# Row -> new value for worksheet "展览" (sheet1)
$exhibitionUpdates = @{
    2  = 41
    4  = 16367
    6  = 20
    8  = 15596
    9  = 68
    10 = 9253
    11 = 473
    12 = 11
    13 = 1027
    14 = 121
    15 = 218
    19 = 91
    20 = 610
    28 = 528
    36 = 365
    37 = 477
    39 = 5679
}

# Row -> new value for worksheet "全部类型" (sheet4)
$allTypesUpdates = @{
    2  = 41
    4  = 16367
    6  = 20
    8  = 15596
    9  = 68
    10 = 9253
    11 = 473
    12 = 11
    13 = 1027
    14 = 121
    15 = 218
    19 = 91
    20 = 610
    28 = 528
    38 = 365
    39 = 477
    41 = 5679
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
